# Swaps upload. Identifying Fixed swaps as well as Bloomberg dividend
# Update avg_long (col U) and avg_short (col V) values on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (XAUUSD)
$ws.Range("U4").Value = -3.5706
$ws.Range("V4").Value = -1.9448
# Row 7 (.A50)
$ws.Range("U7").Value = 64.495
$ws.Range("V7").Value = -98.164
# Row 8 (.AUS200)
$ws.Range("U8").Value = -4.3963
$ws.Range("V8").Value = -6.1616
# Row 9 (.DE30)
$ws.Range("U9").Value = -77.5291
$ws.Range("V9").Value = -154.9163
# Row 10 (.ES35)
$ws.Range("U10").Value = 30.1014
$ws.Range("V10").Value = -53.5347
# Row 11 (.F40)
$ws.Range("U11").Value = -2.0678
$ws.Range("V11").Value = -7.5332
# Row 12 (.HK50)
$ws.Range("U12").Value = 77.551
$ws.Range("V12").Value = -138.6
# Row 13 (.JP225)
$ws.Range("U13").Value = -19.8042
$ws.Range("V13").Value = -21.395
# Row 14 (.STOXX50)
$ws.Range("U14").Value = 0.2718
$ws.Range("V14").Value = -6.5608
# Row 15 (.UK100)
$ws.Range("U15").Value = 0.6367
$ws.Range("V15").Value = -11.812
# Row 17 (.US100)
$ws.Range("U17").Value = -7.9873
$ws.Range("V17").Value = -13.4267
# Row 18 (.US30)
$ws.Range("U18").Value = -13.1026
$ws.Range("V18").Value = -42.1853
# Row 19 (.US500)
$ws.Range("U19").Value = -7.6658
$ws.Range("V19").Value = -57.0889
# Row 20 (.USOil)
$ws.Range("U20").Value = 2.9878
$ws.Range("V20").Value = -23.012
# Row 21 (.XNGUSD)
$ws.Range("U21").Value = -1.958
$ws.Range("V21").Value = 0.3834
# Row 22 (AUDCAD)
$ws.Range("U22").Value = -2.3287
$ws.Range("V22").Value = -0.5929
# Row 23 (AUDCHF)
$ws.Range("U23").Value = 0.3519
$ws.Range("V23").Value = -2.783
# Row 24 (AUDJPY)
$ws.Range("U24").Value = -1.1374
$ws.Range("V24").Value = -1.3365
# Row 25 (AUDNZD)
$ws.Range("U25").Value = -2.3298
$ws.Range("V25").Value = -1.045
# Row 26 (AUDSGD)
$ws.Range("U26").Value = -3.2604
$ws.Range("V26").Value = -1.1363
# Row 27 (AUDUSD)
$ws.Range("U27").Value = -2.8105
$ws.Range("V27").Value = -1.6082
# Row 28 (CADCHF)
$ws.Range("U28").Value = 0.819
$ws.Range("V28").Value = -3.6641
# Row 29 (CADJPY)
$ws.Range("U29").Value = -0.3597
$ws.Range("V29").Value = -2.2407
# Row 30 (CADSGD)
$ws.Range("U30").Value = 0.0477
$ws.Range("V30").Value = -1.6698
# Row 31 (CHFJPY)
$ws.Range("U31").Value = -4.4858
$ws.Range("V31").Value = 0.6534
# Row 32 (CHFSGD)
$ws.Range("U32").Value = -9.906599999999999
$ws.Range("V32").Value = -2.7599
# Row 33 (EURAUD)
$ws.Range("U33").Value = -5.225
$ws.Range("V33").Value = 0.1117
# Row 34 (EURCAD)
$ws.Range("U34").Value = -6.1325
$ws.Range("V34").Value = 0.9765
# Row 35 (EURCHF)
$ws.Range("U35").Value = -0.8943
$ws.Range("V35").Value = -2.5531
# Row 37 (EURGBP)
$ws.Range("U37").Value = -4.0777
$ws.Range("V37").Value = -0.3366
# Row 39 (EURJPY)
$ws.Range("U39").Value = -4.1151
$ws.Range("V39").Value = 0.0511
# Row 42 (EURNZD)
$ws.Range("U42").Value = -6.4878
$ws.Range("V42").Value = 0.7272
# Row 45 (EURSGD)
$ws.Range("U45").Value = -7.7869
$ws.Range("V45").Value = 0.2745
# Row 47 (EURUSD)
$ws.Range("U47").Value = -5.6172
$ws.Range("V47").Value = 0.5976
# Row 48 (GBPAUD)
$ws.Range("U48").Value = -3.432
$ws.Range("V48").Value = -5.291
# Row 49 (GBPCAD)
$ws.Range("U49").Value = -4.8686
$ws.Range("V49").Value = -3.4144
# Row 50 (GBPCHF)
$ws.Range("U50").Value = 0.3645
$ws.Range("V50").Value = -6.8937
# Row 52 (GBPJPY)
$ws.Range("U52").Value = -2.5421
$ws.Range("V52").Value = -4.4506
# Row 54 (GBPNZD)
$ws.Range("U54").Value = -5.3988
$ws.Range("V54").Value = -5.07
# Row 56 (GBPSGD)
$ws.Range("U56").Value = -5.0479
$ws.Range("V56").Value = -3.0767
# Row 57 (GBPUSD)
$ws.Range("U57").Value = -3.9072
$ws.Range("V57").Value = -3.0008
# Row 62 (NZDCAD)
$ws.Range("U62").Value = -1.6357
$ws.Range("V62").Value = -1.034
# Row 63 (NZDCHF)
$ws.Range("U63").Value = 0.5652
$ws.Range("V63").Value = -2.9832
# Row 64 (NZDJPY)
$ws.Range("U64").Value = -0.5775
$ws.Range("V64").Value = -1.6775
# Row 65 (NZDSGD)
$ws.Range("U65").Value = -2.4816
$ws.Range("V65").Value = -1.6522
# Row 75 (USDJPY)
$ws.Range("U75").Value = -1.4196
$ws.Range("V75").Value = -3.6456
